{"js": "// Helper: find the index of the first paragraph (from `fromIndex` onward)\n// whose text satisfies `predicate`.\nfunction findParagraphIndex(items, predicate, fromIndex) {\n  const start = fromIndex || 0;\n  for (let i = start; i < items.length; i++) {\n    if (predicate(items[i].text)) {\n      return i;\n    }\n  }\n  return -1;\n}\n\n// Helper: replace the first substring match of `find` inside a single\n// paragraph's own range with `replacement`, leaving the rest of the\n// paragraph's runs untouched.\nasync function replaceInParagraph(paragraph, find, replacement) {\n  const range = paragraph.getRange();\n  const results = range.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Locate the six existing \"queries\" paragraphs -------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet items = paragraphs.items;\n\nconst idx1 = findParagraphIndex(items, (t) => t.indexOf(\"Select percentage of population confirmed\") !== -1);\nconst idx2 = findParagraphIndex(items, (t) => t.indexOf(\"Select percentage of population recovered\") !== -1, idx1 + 1);\nconst idx3 = findParagraphIndex(items, (t) => t.indexOf(\"Select percentage of population dead\") !== -1, idx2 + 1);\nconst idx4 = findParagraphIndex(items, (t) => t.indexOf(\"Cumulative average confirmed\") !== -1, idx3 + 1);\nconst idx5 = findParagraphIndex(items, (t) => t.indexOf(\"Cumulative average death\") !== -1, idx4 + 1);\nconst idx6 = findParagraphIndex(items, (t) => t.indexOf(\"Cumulative average recovers\") !== -1, idx5 + 1);\n\n// 1) Paragraph \"1 \\u2013 Select percentage of population confirmed.\" gains\n//    \", dead, and recovered\" before the final period.\nawait replaceInParagraph(\n  items[idx1],\n  \"Select percentage of population confirmed.\",\n  \"Select percentage of population confirmed, dead, and recovered.\"\n);\n\n// 2) Paragraph \"4 \\u2013 Cumulative average confirmed.\" becomes item \"2\" and\n//    gains \", deaths, and recovers\" before the final period.\nawait replaceInParagraph(items[idx4], \"4 \\u2013\", \"2 \\u2013\");\nawait replaceInParagraph(\n  items[idx4],\n  \"Cumulative average confirmed.\",\n  \"Cumulative average confirmed, deaths, and recovers.\"\n);\n\n// Remove the now-superseded paragraphs (old items 2, 3, 5, 6).\nitems[idx2].delete();\nitems[idx3].delete();\nitems[idx5].delete();\nitems[idx6].delete();\nawait context.sync();\n\n// Re-load to get a fresh, valid reference to the (renumbered) item \"2\"\n// paragraph so we can insert the four brand-new questions after it.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\nitems = paragraphs2.items;\n\nconst item2Idx = findParagraphIndex(items, (t) => t.indexOf(\"Cumulative average confirmed, deaths, and recovers\") !== -1);\n\nconst p3 = items[item2Idx].insertParagraph(\n  \"3 - When was the peak of mortality rate of the pandemic?\",\n  Word.InsertLocation.after\n);\nconst p4 = p3.insertParagraph(\n  \"4 - What are the top 10 countries with greatest amount of deaths?\",\n  Word.InsertLocation.after\n);\nconst p5 = p4.insertParagraph(\n  \"5 - What are the correlation between deaths and population?\",\n  Word.InsertLocation.after\n);\np5.insertParagraph(\n  \"6 - What is the average recovered rate by countries?\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// --- Add the new \"Stretches\" bullet after \"Use AWS or Azure.\" -------------\nconst paragraphs3 = context.document.body.paragraphs;\nparagraphs3.load(\"items/text\");\nawait context.sync();\n\nconst awsIdx = findParagraphIndex(paragraphs3.items, (t) => t.indexOf(\"AWS or Azure\") !== -1);\n\nparagraphs3.items[awsIdx].insertParagraph(\n  \"Find relations with other kind of Data, i.e., GDP.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex {\n    param($doc, [string]$needle, [int]$startAfter = 0)\n    $i = 0\n    foreach ($p in $doc.Paragraphs) {\n        $i = $i + 1\n        if ($i -gt $startAfter -and $p.Range.Text -like \"*$needle*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Replace-InParagraph {\n    param($doc, [int]$paraIndex, [string]$find, [string]$replacement)\n    $p = $doc.Paragraphs($paraIndex)\n    $rng = $p.Range.Duplicate\n    $rng.Find.Text = $find\n    $rng.Find.Execute() | Out-Null\n    $rng.Text = $replacement\n}\n\n# --- Locate the six existing \"queries\" paragraphs --------------------------\n$idx1 = Find-ParagraphIndex $d \"Select percentage of population confirmed\"\n$idx2 = Find-ParagraphIndex $d \"Select percentage of population recovered\" $idx1\n$idx3 = Find-ParagraphIndex $d \"Select percentage of population dead\" $idx2\n$idx4 = Find-ParagraphIndex $d \"Cumulative average confirmed\" $idx3\n$idx5 = Find-ParagraphIndex $d \"Cumulative average death\" $idx4\n$idx6 = Find-ParagraphIndex $d \"Cumulative average recovers\" $idx5\n\n# 1) Paragraph \"1 - Select percentage of population confirmed.\" gains\n#    \", dead, and recovered\" before the final period.\nReplace-InParagraph $d $idx1 \"Select percentage of population confirmed.\" \"Select percentage of population confirmed, dead, and recovered.\"\n\n# 2) Paragraph \"4 - Cumulative average confirmed.\" becomes item \"2\" and\n#    gains \", deaths, and recovers\" before the final period.\nReplace-InParagraph $d $idx4 \"4 \u2013\" \"2 \u2013\"\nReplace-InParagraph $d $idx4 \"Cumulative average confirmed.\" \"Cumulative average confirmed, deaths, and recovers.\"\n\n# Remove the now-superseded paragraphs (old items 2, 3, 5, 6), deleting from\n# the highest index down so earlier indices stay valid.\n$toDelete = @($idx2, $idx3, $idx5, $idx6) | Sort-Object -Descending\nforeach ($i in $toDelete) {\n    $d.Paragraphs($i).Range.Delete()\n}\n\n# Insert the four brand-new question paragraphs right after the\n# (renumbered) item \"2\" paragraph.\n$item2Idx = Find-ParagraphIndex $d \"Cumulative average confirmed, deaths, and recovers\"\n$p = $d.Paragraphs($item2Idx)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs($item2Idx + 1).Range.Text = \"3 - When was the peak of mortality rate of the pandemic?\"\n\n$p = $d.Paragraphs($item2Idx + 1)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs($item2Idx + 2).Range.Text = \"4 - What are the top 10 countries with greatest amount of deaths?\"\n\n$p = $d.Paragraphs($item2Idx + 2)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs($item2Idx + 3).Range.Text = \"5 - What are the correlation between deaths and population?\"\n\n$p = $d.Paragraphs($item2Idx + 3)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs($item2Idx + 4).Range.Text = \"6 - What is the average recovered rate by countries?\"\n\n# --- Add the new \"Stretches\" bullet after \"Use AWS or Azure.\" --------------\n$awsIdx = Find-ParagraphIndex $d \"AWS or Azure\"\n$p = $d.Paragraphs($awsIdx)\n$p.Range.InsertParagraphAfter()\n$d.Paragraphs($awsIdx + 1).Range.Text = \"Find relations with other kind of Data, i.e., GDP.\"\n"}
